$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").ClearContents()
$ws.Range("A1").Value = "working on little feature"
$ws.Range("A1").Select()
